# Commit: Wed, Jul 08, 2020 10:05:18 PM
#
# The canonical diff does two things:
#   1. Swaps the two embedded theme parts' color palettes: the Slide
#      Master's theme (ppt/theme/theme1.xml -- the "Integral"/"Red Violet"
#      theme) ends up with the default "Office Theme" color palette that
#      used to live in ppt/theme/theme2.xml (the Notes Master's theme),
#      and vice versa.
#   2. Re-points the three "Component three" tables that used the old
#      custom table-style GUID at a new built-in table-style GUID.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1. Theme colors
# ---------------------------------------------------------------------
# Target palette (the former "Office Theme" colors from theme2.xml),
# expressed as COM RGB() integers (0xBBGGRR) in clrScheme order:
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink.
$officeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$colorScheme = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeColors[$i - 1]
}

# ---------------------------------------------------------------------
# 2. Table styles: old custom style GUID -> new built-in style GUID
# ---------------------------------------------------------------------
$oldStyleId = "{3363C96D-39F7-4E5B-A7BB-566259EC3A2B}"
$newStyleId = "{618C6CA9-A5C4-4A07-880D-C309BF49E5E0}"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $shp = $sl.Shapes.Item($shi)
        if ($shp.HasTable) {
            $tbl = $shp.Table
            if ($tbl.Style -eq $oldStyleId) {
                $tbl.ApplyStyle($newStyleId)
            }
        }
    }
}
